# Auto-generated Excel COM-interop script to append new sensor-log rows
# to the PIR, Humidity, and Temperature sheets, matching the upstream diff.

$wb = $excel.ActiveWorkbook

# --- PIR sheet ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A447").Value = "'2026-02-06"
$ws.Range("A447").Style = "Normal"
$ws.Range("B447").Value = "10:17:03"
$ws.Range("C447").Value = "10:00"
$ws.Range("D447").Value = "Bathroom"
$ws.Range("E447").Value = "No Motion"
$ws.Range("F447").Value = "Inactive"
$ws.Range("A448").Value = "'2026-02-06"
$ws.Range("A448").Style = "Normal"
$ws.Range("B448").Value = "10:17:04"
$ws.Range("C448").Value = "10:00"
$ws.Range("D448").Value = "Bathroom"
$ws.Range("E448").Value = "No Motion"
$ws.Range("F448").Value = "Inactive"
$ws.Range("A449").Value = "'2026-02-06"
$ws.Range("A449").Style = "Normal"
$ws.Range("B449").Value = "10:17:09"
$ws.Range("C449").Value = "10:00"
$ws.Range("D449").Value = "Bathroom"
$ws.Range("E449").Value = "No Motion"
$ws.Range("F449").Value = "Inactive"
$ws.Range("A450").Value = "'2026-02-06"
$ws.Range("A450").Style = "Normal"
$ws.Range("B450").Value = "10:17:14"
$ws.Range("C450").Value = "10:00"
$ws.Range("D450").Value = "Bathroom"
$ws.Range("E450").Value = "No Motion"
$ws.Range("F450").Value = "Inactive"
$ws.Range("A451").Value = "'2026-02-06"
$ws.Range("A451").Style = "Normal"
$ws.Range("B451").Value = "10:17:19"
$ws.Range("C451").Value = "10:00"
$ws.Range("D451").Value = "Bathroom"
$ws.Range("E451").Value = "No Motion"
$ws.Range("F451").Value = "Inactive"
$ws.Range("A452").Value = "'2026-02-06"
$ws.Range("A452").Style = "Normal"
$ws.Range("B452").Value = "10:17:24"
$ws.Range("C452").Value = "10:00"
$ws.Range("D452").Value = "Bathroom"
$ws.Range("E452").Value = "No Motion"
$ws.Range("F452").Value = "Inactive"
$ws.Range("A453").Value = "'2026-02-06"
$ws.Range("A453").Style = "Normal"
$ws.Range("B453").Value = "10:17:29"
$ws.Range("C453").Value = "10:00"
$ws.Range("D453").Value = "Bathroom"
$ws.Range("E453").Value = "No Motion"
$ws.Range("F453").Value = "Inactive"
$ws.Range("A454").Value = "'2026-02-06"
$ws.Range("A454").Style = "Normal"
$ws.Range("B454").Value = "10:17:34"
$ws.Range("C454").Value = "10:00"
$ws.Range("D454").Value = "Bathroom"
$ws.Range("E454").Value = "No Motion"
$ws.Range("F454").Value = "Inactive"
$ws.Range("A455").Value = "'2026-02-06"
$ws.Range("A455").Style = "Normal"
$ws.Range("B455").Value = "10:17:39"
$ws.Range("C455").Value = "10:00"
$ws.Range("D455").Value = "Bathroom"
$ws.Range("E455").Value = "No Motion"
$ws.Range("F455").Value = "Inactive"
$ws.Range("A456").Value = "'2026-02-06"
$ws.Range("A456").Style = "Normal"
$ws.Range("B456").Value = "10:17:44"
$ws.Range("C456").Value = "10:00"
$ws.Range("D456").Value = "Bathroom"
$ws.Range("E456").Value = "No Motion"
$ws.Range("F456").Value = "Inactive"
$ws.Range("A457").Value = "'2026-02-06"
$ws.Range("A457").Style = "Normal"
$ws.Range("B457").Value = "10:17:49"
$ws.Range("C457").Value = "10:00"
$ws.Range("D457").Value = "Bathroom"
$ws.Range("E457").Value = "No Motion"
$ws.Range("F457").Value = "Inactive"
$ws.Range("A458").Value = "'2026-02-06"
$ws.Range("A458").Style = "Normal"
$ws.Range("B458").Value = "10:17:54"
$ws.Range("C458").Value = "10:00"
$ws.Range("D458").Value = "Bathroom"
$ws.Range("E458").Value = "No Motion"
$ws.Range("F458").Value = "Inactive"
$ws.Range("A459").Value = "'2026-02-06"
$ws.Range("A459").Style = "Normal"
$ws.Range("B459").Value = "10:17:59"
$ws.Range("C459").Value = "10:00"
$ws.Range("D459").Value = "Bathroom"
$ws.Range("E459").Value = "No Motion"
$ws.Range("F459").Value = "Inactive"

# --- Humidity sheet ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A309").Value = "'2026-02-06"
$ws.Range("A309").Style = "Normal"
$ws.Range("B309").Value = "10:17:05"
$ws.Range("C309").Value = "10:00"
$ws.Range("D309").Value = "Bathroom"
$ws.Range("E309").Value = "'67.1%"
$ws.Range("E309").Style = "Normal"
$ws.Range("F309").Value = "Active"
$ws.Range("A310").Value = "'2026-02-06"
$ws.Range("A310").Style = "Normal"
$ws.Range("B310").Value = "10:17:10"
$ws.Range("C310").Value = "10:00"
$ws.Range("D310").Value = "Bathroom"
$ws.Range("E310").Value = "'67.4%"
$ws.Range("E310").Style = "Normal"
$ws.Range("F310").Value = "Active"
$ws.Range("A311").Value = "'2026-02-06"
$ws.Range("A311").Style = "Normal"
$ws.Range("B311").Value = "10:17:15"
$ws.Range("C311").Value = "10:00"
$ws.Range("D311").Value = "Bathroom"
$ws.Range("E311").Value = "'68.4%"
$ws.Range("E311").Style = "Normal"
$ws.Range("F311").Value = "Active"
$ws.Range("A312").Value = "'2026-02-06"
$ws.Range("A312").Style = "Normal"
$ws.Range("B312").Value = "10:17:25"
$ws.Range("C312").Value = "10:00"
$ws.Range("D312").Value = "Bathroom"
$ws.Range("E312").Value = "'68.4%"
$ws.Range("E312").Style = "Normal"
$ws.Range("F312").Value = "Active"
$ws.Range("A313").Value = "'2026-02-06"
$ws.Range("A313").Style = "Normal"
$ws.Range("B313").Value = "10:17:45"
$ws.Range("C313").Value = "10:00"
$ws.Range("D313").Value = "Bathroom"
$ws.Range("E313").Value = "'68.6%"
$ws.Range("E313").Style = "Normal"
$ws.Range("F313").Value = "Active"
$ws.Range("A314").Value = "'2026-02-06"
$ws.Range("A314").Style = "Normal"
$ws.Range("B314").Value = "10:17:55"
$ws.Range("C314").Value = "10:00"
$ws.Range("D314").Value = "Bathroom"
$ws.Range("E314").Value = "'68.8%"
$ws.Range("E314").Style = "Normal"
$ws.Range("F314").Value = "Active"

# --- Temperature sheet ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A308").Value = "'2026-02-06"
$ws.Range("A308").Style = "Normal"
$ws.Range("B308").Value = "10:17:01"
$ws.Range("C308").Value = "10:00"
$ws.Range("D308").Value = "Bathroom"
$ws.Range("E308").Value = "28.3C"
$ws.Range("F308").Value = "Active"
$ws.Range("A309").Value = "'2026-02-06"
$ws.Range("A309").Style = "Normal"
$ws.Range("B309").Value = "10:17:06"
$ws.Range("C309").Value = "10:00"
$ws.Range("D309").Value = "Bathroom"
$ws.Range("E309").Value = "28.2C"
$ws.Range("F309").Value = "Active"
$ws.Range("A310").Value = "'2026-02-06"
$ws.Range("A310").Style = "Normal"
$ws.Range("B310").Value = "10:17:11"
$ws.Range("C310").Value = "10:00"
$ws.Range("D310").Value = "Bathroom"
$ws.Range("E310").Value = "28.1C"
$ws.Range("F310").Value = "Active"
$ws.Range("A311").Value = "'2026-02-06"
$ws.Range("A311").Style = "Normal"
$ws.Range("B311").Value = "10:17:16"
$ws.Range("C311").Value = "10:00"
$ws.Range("D311").Value = "Bathroom"
$ws.Range("E311").Value = "28.2C"
$ws.Range("F311").Value = "Active"
$ws.Range("A312").Value = "'2026-02-06"
$ws.Range("A312").Style = "Normal"
$ws.Range("B312").Value = "10:17:26"
$ws.Range("C312").Value = "10:00"
$ws.Range("D312").Value = "Bathroom"
$ws.Range("E312").Value = "28.2C"
$ws.Range("F312").Value = "Active"
$ws.Range("A313").Value = "'2026-02-06"
$ws.Range("A313").Style = "Normal"
$ws.Range("B313").Value = "10:17:46"
$ws.Range("C313").Value = "10:00"
$ws.Range("D313").Value = "Bathroom"
$ws.Range("E313").Value = "28.2C"
$ws.Range("F313").Value = "Active"
$ws.Range("A314").Value = "'2026-02-06"
$ws.Range("A314").Style = "Normal"
$ws.Range("B314").Value = "10:17:56"
$ws.Range("C314").Value = "10:00"
$ws.Range("D314").Value = "Bathroom"
$ws.Range("E314").Value = "28.2C"
$ws.Range("F314").Value = "Active"
